# Daily attendance processing - 2026-01-30 15:49:04
#
# The "Recorded By" column (G) lists session recorders as a
# comma-separated string. For sessions recorded by both the automated
# System and dnasr281@gmail.com, the order of the two names is being
# swapped from "System, dnasr281@gmail.com" to
# "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$searchRange = $ws.Columns.Item(7)  # Column G - "Recorded By"

$first = $searchRange.Find($oldValue)
if ($first -ne $null) {
    $firstAddress = $first.Address()
    $current = $first
    do {
        $current.Value = $newValue
        $current = $searchRange.FindNext($current)
    } while ($current -ne $null -and $current.Address() -ne $firstAddress)
}
